# Update the "Colour Code" gradient (column B) on the active sheet.
# The palette was regenerated with one additional shade so the 54 rows
# (2-55) now span a slightly denser blue gradient, and a colour value is
# now also populated for row 55 (previously blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colours = @{
    2  = "#fcfdff"
    3  = "#f9faff"
    4  = "#f6f8ff"
    5  = "#f3f6ff"
    6  = "#f0f4ff"
    7  = "#edf2ff"
    8  = "#ebefff"
    9  = "#e8edff"
    10 = "#e5ebff"
    11 = "#e2e9ff"
    12 = "#dfe6ff"
    13 = "#dce4ff"
    14 = "#d8e2ff"
    15 = "#d5e0ff"
    16 = "#d2deff"
    17 = "#cfdbff"
    18 = "#ccd9ff"
    19 = "#c9d7ff"
    20 = "#c6d5ff"
    21 = "#c3d3ff"
    22 = "#bfd1ff"
    23 = "#bcceff"
    24 = "#b9ccff"
    25 = "#b5caff"
    26 = "#b2c8ff"
    27 = "#afc6ff"
    28 = "#abc4ff"
    29 = "#a8c2ff"
    30 = "#a4c0ff"
    31 = "#a1bdff"
    32 = "#9dbbff"
    33 = "#9ab9ff"
    34 = "#96b7ff"
    35 = "#92b5ff"
    36 = "#8eb3ff"
    37 = "#8ab1ff"
    38 = "#86afff"
    39 = "#82adff"
    40 = "#7eabff"
    41 = "#7aa9ff"
    42 = "#75a7ff"
    43 = "#71a5ff"
    44 = "#6ca3ff"
    45 = "#67a1ff"
    46 = "#629fff"
    47 = "#5c9dff"
    48 = "#569bff"
    49 = "#5099ff"
    50 = "#4a97ff"
    51 = "#4295ff"
    52 = "#3a93ff"
    53 = "#3091ff"
    54 = "#248fff"
    55 = "#118dff"
}

foreach ($row in $colours.Keys) {
    $ws.Range("B$row").Value = $colours[$row]
}
